$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 1) "ROOT" sheet: sort the data rows (A2:H25) ascending by column A,
#    matching the Data > Sort operation performed by the author.
# ----------------------------------------------------------------------
$wsRoot = $wb.Worksheets.Item("ROOT")
$wsRoot.Activate()

$sortObj = $wsRoot.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($wsRoot.Range("A1"))
$sortObj.SetRange($wsRoot.Range("A1:H25"))
$sortObj.Header = 1
$sortObj.Apply()

# The sort engine leaves a handful of empty-but-still-formatted cells
# behind at the old location of the row that used to carry custom
# formatting (old row 23, now relocated to row 4). Clear them so the
# sheet doesn't retain stray formatting-only cells.
$wsRoot.Range("E23:H23").Clear()

# Restore the active cell on this sheet to B14 (matches the author's
# final selection after sorting).
$wsRoot.Range("B14").Select()

# ----------------------------------------------------------------------
# 2) "Tables" sheet: the cable landing points row used to reference the
#    (non-existent) "project.farm" table via "fk_site_id"; point it at
#    "project.site" via "site_name" instead.
# ----------------------------------------------------------------------
$wsTables = $wb.Worksheets.Item("Tables")
$wsTables.Activate()

$wsTables.Range("B7").Value = "project.site"
$wsTables.Range("C7").Value = "site_name"

# Update the active cell on this sheet to C7.
$wsTables.Range("C7").Select()
